# Plano de Ação BeCold - sprint week update
# 1) Sort the "Dados" table by column D (Semana) ascending, matching the
#    table's sortState (B7:H21 sorted by D6:D21).
# 2) Update Responsável/Atribuído a (and a couple of status/date) cells on
#    several rows that received assignees this sprint.
# 3) Refresh the sheet view (scroll position / active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados do plano de Ação")
$ws.Activate()

$lo = $ws.ListObjects.Item("Dados")

$lo.Sort.SortFields.Clear()
$lo.Sort.SortFields.Add($ws.Range("D6:D21"))
$lo.Sort.Header = 1
$lo.Sort.Apply()

# "Script Banco de Dados" -> Concluído, assigned to Thalita, Guilherme
$ws.Range("C15").Value = "Concluído"
$ws.Range("F15").Value = "Thalita, Guilherme"

# "Página Inicial Estática" -> Thalita, assigned to Thalita, Guilherme; dates moved up
$ws.Range("E16").Value = "Thalita"
$ws.Range("F16").Value = "Thalita, Guilherme"
$ws.Range("G16").Formula = "=DATE(YEAR(TODAY()),10,4)"
$ws.Range("H16").Formula = "=DATE(YEAR(TODAY()),10,11)"

# "Página Login Estática" -> Thalita, assigned to Yuri
$ws.Range("E17").Value = "Thalita"
$ws.Range("F17").Value = "Yuri"

# "Página Cadastro Estática" -> Thalita, assigned to Paulo, Kaiqui; end date moved up
$ws.Range("E18").Value = "Thalita"
$ws.Range("F18").Value = "Paulo, Kaiqui"
$ws.Range("H18").Formula = "=DATE(YEAR(TODAY()),10,11)"

# "Página Dashboard Estática" -> Thalita
$ws.Range("E19").Value = "Thalita"

$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("F17").Select()
